$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 3-9 with new data (each row effectively shifted up with some values changed),
# and remove the old row 10 (Urine LAM) since the table now only spans A1:C9.

$ws.Range("A3").Value = "VOC Model (TPP Optimized)"
$ws.Range("B3").Value = 84.40000000000001
$ws.Range("C3").Value = 80

$ws.Range("A4").Value = "VOC Model (Sensitivity Optimized)"
$ws.Range("B4").Value = 93.8
$ws.Range("C4").Value = 26.7

$ws.Range("A5").Value = "Combined Sensitivity Optimized with CXR (Parallel)"
$ws.Range("B5").Value = 96.3
$ws.Range("C5").Value = 65.59999999999999

$ws.Range("A6").Value = "CXR + CAD"
$ws.Range("B6").Value = 90
$ws.Range("C6").Value = 74

$ws.Range("A7").Value = "CXR"
$ws.Range("B7").Value = 76
$ws.Range("C7").Value = 82

$ws.Range("A8").Value = "CRP"
$ws.Range("B8").Value = 84
$ws.Range("C8").Value = 61

$ws.Range("A9").Value = "Urine LAM"
$ws.Range("B9").Value = 42
$ws.Range("C9").Value = 99

# Remove the now-obsolete row 10 entirely so the used range shrinks to A1:C9
$ws.Rows.Item(10).Delete()
